# "edited my user story"
# Append a new "As a Delivery Driver" user story row for printing an
# invoice, then leave the new cell selected (matching the cursor position
# recorded by Excel when the author saved the file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "As a Delivery Driver"
$ws.Range("B18").Value = "I want to print invoice"

$ws.Range("B18").Select()
